# Update cryptocurrency Price / Volume(1h) figures on Sheet1 to reflect
# the latest scrape, matching the GitHub Actions commit on
# Sun Jun  2 23:51:13 UTC 2024.
#
# Values are written as literal text (not auto-converted numbers) to
# match the original inline-string cells produced by the scraper, using
# a leading apostrophe to force text entry and then resetting the style
# back to Normal so no stray "quote prefix" number format lingers on the
# cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")


$ws.Range("D2").Value = '67.846.42'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '3.784.58'
$ws.Range("E3").Value = '  -0.82%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'603.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("D6").Value = "'163.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.52%  '
$ws.Range("D7").Value = '3.783.20'
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("E8").Value = '  -0.15%  '
$ws.Range("E9").Value = '  -0.91%  '
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("D11").Value = "'0.447"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.63%  '
$ws.Range("D12").Value = "'6.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.17%  '
$ws.Range("E13").Value = '  -2.38%  '
$ws.Range("D14").Value = "'35.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").Value = '4.417.31'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").Value = '3.780.41'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").Value = '67.828.31'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = "'18.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.27%  '
$ws.Range("E19").Value = '  +1.90%  '
$ws.Range("D20").Value = "'7.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").Value = "'458.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("D22").Value = "'9.44"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.50%  '
$ws.Range("D23").Value = "'0.691"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.12%  '
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").Value = "'83.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'11.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.62%  '
$ws.Range("E27").Value = '  -1.04%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").Value = "'9.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("D30").Value = '3.929.34'
$ws.Range("E30").Value = '  -0.90%  '
$ws.Range("E31").Value = '  -6.48%  '
$ws.Range("D32").Value = "'7.20"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.18%  '
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("D34").Value = "'28.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("D37").Value = "'0.0993"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.58%  '
$ws.Range("E38").Value = '  +7.37%  '
$ws.Range("E39").Value = '  +0.30%  '
$ws.Range("D40").Value = "'3.20"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.27%  '
$ws.Range("E41").Value = '  -1.90%  '
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("E43").Value = '  +0.16%  '
$ws.Range("D44").Value = "'43.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.79%  '
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("D46").Value = "'152.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.60%  '
$ws.Range("E47").Value = '  -1.79%  '
$ws.Range("E48").Value = '  -1.92%  '
$ws.Range("E49").Value = '  -0.26%  '
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").Value = "'26.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.11%  '
